# Monitoreo a las actividades del 29 de abril al 6 de mayo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("No Conformidades")

# --- Row 23 : ID 20 ---
$ws.Range("B23").Value = "No todas las tareas se encontraban cronometradas"
$ws.Range("C23").Value = "Ventas"
$ws.Range("D23").Value = 42557
$ws.Range("E23").Value = 42557
$ws.Range("F23").Value = "Cerrada"
$ws.Range("G23").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."

# --- Row 24 : ID 21 ---
$ws.Range("B24").Value = "Las tareas no fueron completadas en su tiempo"
$ws.Range("C24").Value = "Compras"
$ws.Range("D24").Value = 42557
$ws.Range("E24").Value = 42557
$ws.Range("F24").Value = "Cerrada"
$ws.Range("G24").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."

# Row height for the two newly-filled rows
$ws.Range("A23:G24").RowHeight = 75

# Styles for D/E/G on rows 25 and 26 shift to the "filled" look (date format /
# wrap-left) even though those cells stay empty.
$ws.Range("D25:E26").NumberFormat = "m/d/yyyy"
$ws.Range("G25:G26").WrapText = $true
$ws.Range("G25:G26").HorizontalAlignment = -4131

# --- View state ---
$ws.Range("G26:G28").Select()
$excel.ActiveWindow.ScrollRow = 23
